$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin-list rows: prices (col D, forced as literal text with a leading
# apostrophe so Excel keeps storing them as text rather than coercing to a
# float/losing trailing zeros), plus a handful of Coin/Link/Volume(1h) cells
# that were rearranged by the upstream data refresh.

$ws.Range("D2").Value = '''250.98'
$ws.Range("D3").Value = '''22.88'
$ws.Range("D4").Value = '''5.421'
$ws.Range("D5").Value = '''0.05669'
$ws.Range("D7").Value = '''6.384'
$ws.Range("D8").Value = '''0.8127'
$ws.Range("D9").Value = '''0.9247'
$ws.Range("D11").Value = '''0.07455'
$ws.Range("D12").Value = '''0.03194'
$ws.Range("D13").Value = '''0.03065'
$ws.Range("D14").Value = '''0.09356'
$ws.Range("D15").Value = '''3.720'
$ws.Range("D16").Value = '''0.001602'
$ws.Range("D17").Value = '''0.04755'
$ws.Range("D18").Value = '''0.0005793'
$ws.Range("E18").Value = '17OneONEWorstin24h'
$ws.Range("D19").Value = '''0.006369'
$ws.Range("D20").Value = '''0.005044'
$ws.Range("D21").Value = '''0.001028'
$ws.Range("D22").Value = '''0.0001500'
$ws.Range("D23").Value = '''3.707'
$ws.Range("D24").Value = '''2.174'
$ws.Range("D25").Value = '''0.3303'
$ws.Range("D26").Value = '''0.1309'
$ws.Range("D28").Value = '''0.0003002'
$ws.Range("D40").Value = '''0.04022'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = '''0.006779'
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = '''0.1070'
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = '''0.002710'
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("D44").Value = '''0.007550'
$ws.Range("D45").Value = '''0.00005807'
$ws.Range("D47").Value = '''0.5003'
